$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values per repull/mean recalculation
$values = @{
    2  = -1
    3  = 3
    4  = -5
    5  = -2
    6  = 7
    7  = 2
    8  = -2
    9  = 2
    11 = 2
    12 = 2
    13 = -6
    14 = -1
    15 = 1
    16 = -1
    17 = 1
    18 = 1
    20 = -4
    21 = 3
    23 = 1
    24 = 1
    25 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
